# Apply updated cryptocurrency price (D) and 1h volume change (E) figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.387.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.254.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.00%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.47%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.245.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.98%  "

$ws.Range("E10").Value = "  -11.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.578"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.77"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "629.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.786.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.420.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.116"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.260.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.896"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.94"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.96%  "

$ws.Range("E34").Value = "  -5.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.713.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "520.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0727"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.130"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.333"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.23"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0410"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.44%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("E51").Value = "  +0.10%  "
